$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates for rows 2-27, 40 (text-valued numbers) ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "252.60"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.00"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.541"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05695"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.456"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8072"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.041"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1433"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07310"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03143"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02937"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09277"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001677"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.205"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04777"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005807"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006452"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005071"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001052"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001500"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.988"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.381"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.114"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003098"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04138"

# --- Row 17: Volume(1h) label change ---
$ws.Range("E17").Value = "16OneONEWorstin24h"

# --- Row 41: BKEXToken -> KickToken ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006890"
$ws.Range("E41").Value = "40KickTokenKICK"

# --- Row 42: CEJI price update ---
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003500"

# --- Row 43: KickToken -> BKEXToken ---
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1046"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- Price (column D) updates for rows 44-48 ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009543"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005646"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7848"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.01700"
